# AFDP-1050 - Implement Document level security - configuration file updates
#
# Adds a new "Folder" access-control rule row (Folder - deny no access) to
# the rule table on Sheet1, directly below the existing "Folder - default
# public access" row (row 39), mirroring its layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of the last existing rule row (row 39, the
# "Folder" rule) down into the new row 40, cell by cell, so the new row
# keeps the same borders/fill/font/number-format as the rest of the table.
$ws.Range("A39").Copy($ws.Range("A40"))
$ws.Range("B39").Copy($ws.Range("B40"))
$ws.Range("C39").Copy($ws.Range("C40"))
$ws.Range("D39").Copy($ws.Range("D40"))
$ws.Range("E39").Copy($ws.Range("E40"))
$ws.Range("F39").Copy($ws.Range("F40"))
$ws.Range("G39").Copy($ws.Range("G40"))

# New rule: Folder - deny no access
$ws.Range("B40").Value = "Folder –deny no access"
$ws.Range("C40").Value = "FOLDER"
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = ""
$ws.Range("G40").Value = "mandatory deny read to No Access"

$ws.Rows.Item(40).RowHeight = 13.8

# Match the saved selection state from the authored change.
$ws.Range("B40").Select() | Out-Null
